# Applies the "Updated symbol list on Wed Jan  4 09:51:44 UTC 2023 with GitHub Actions"
# data refresh to the cryptos worksheet: refreshed prices / 1h volumes, plus a
# handful of rows whose coin (and therefore rank order) changed.
#
# Values in columns D (Price) and E (Volume(1h)) are stored as literal text in
# the source workbook (e.g. "255.33", "3.73%") rather than numbers, so we
# prefix numeric-looking values with a leading apostrophe to force Excel to
# keep them as text instead of auto-converting to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "'255.33" }
    @{ Cell = "E2"; Value = "'3.73%" }
    @{ Cell = "D3"; Value = "'28.24" }
    @{ Cell = "E3"; Value = "'-5.24%" }
    @{ Cell = "D4"; Value = "'5.247" }
    @{ Cell = "E4"; Value = "'1.77%" }
    @{ Cell = "D5"; Value = "'0.05851" }
    @{ Cell = "E5"; Value = "'1.57%" }
    @{ Cell = "D6"; Value = "'6.713" }
    @{ Cell = "E6"; Value = "'0.79%" }
    @{ Cell = "B7"; Value = "MXToken" }
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" }
    @{ Cell = "D7"; Value = "'0.8667" }
    @{ Cell = "E7"; Value = "'1.88%" }
    @{ Cell = "B8"; Value = "FTXToken" }
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" }
    @{ Cell = "D8"; Value = "'1.050" }
    @{ Cell = "E8"; Value = "'22.70%" }
    @{ Cell = "B9"; Value = "WazirX" }
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" }
    @{ Cell = "D9"; Value = "'0.1411" }
    @{ Cell = "E9"; Value = "'1.47%" }
    @{ Cell = "B10"; Value = "MandalaExchangeToken" }
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" }
    @{ Cell = "D10"; Value = "'0.07154" }
    @{ Cell = "E10"; Value = "'1.01%" }
    @{ Cell = "B11"; Value = "BitrueCoin" }
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" }
    @{ Cell = "D11"; Value = "'0.03187" }
    @{ Cell = "E11"; Value = "'-1.82%" }
    @{ Cell = "B12"; Value = "BitMartToken" }
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" }
    @{ Cell = "D12"; Value = "'0.09235" }
    @{ Cell = "E12"; Value = "'-1.48%" }
    @{ Cell = "B13"; Value = "BitForexToken" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" }
    @{ Cell = "D13"; Value = "'0.001538" }
    @{ Cell = "E13"; Value = "'0.13%" }
    @{ Cell = "B14"; Value = "One" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" }
    @{ Cell = "D14"; Value = "'0.0006061" }
    @{ Cell = "E14"; Value = "'-94.07%" }
    @{ Cell = "D15"; Value = "'0.005809" }
    @{ Cell = "E15"; Value = "'-2.11%" }
    @{ Cell = "D16"; Value = "'3.499" }
    @{ Cell = "E16"; Value = "'-0.64%" }
    @{ Cell = "B17"; Value = "GateToken" }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" }
    @{ Cell = "D17"; Value = "'3.229" }
    @{ Cell = "E17"; Value = "'-0.29%" }
    @{ Cell = "B18"; Value = "BTSEToken" }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" }
    @{ Cell = "D18"; Value = "'2.223" }
    @{ Cell = "E18"; Value = "'0.85%" }
    @{ Cell = "D19"; Value = "'0.3182" }
    @{ Cell = "E19"; Value = "'0.46%" }
    @{ Cell = "E20"; Value = "'3.30%" }
    @{ Cell = "E21"; Value = "'0.42%" }
    @{ Cell = "D22"; Value = "'3.531" }
    @{ Cell = "E22"; Value = "'0.80%" }
    @{ Cell = "D23"; Value = "'0.04149" }
    @{ Cell = "E23"; Value = "'0.86%" }
    @{ Cell = "E24"; Value = "'-4.40%" }
    @{ Cell = "D25"; Value = "'0.001229" }
    @{ Cell = "E25"; Value = "'0.11%" }
    @{ Cell = "D26"; Value = "'0.004816" }
    @{ Cell = "E26"; Value = "'16.03%" }
    @{ Cell = "E27"; Value = "'0.05%" }
    @{ Cell = "D28"; Value = "'0.0001466" }
    @{ Cell = "E28"; Value = "'1.29%" }
    @{ Cell = "D40"; Value = "'0.03802" }
    @{ Cell = "E40"; Value = "'1.45%" }
    @{ Cell = "B41"; Value = "BKEXToken" }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk" }
    @{ Cell = "D41"; Value = "'0.1103" }
    @{ Cell = "E41"; Value = "'2.84%" }
    @{ Cell = "B42"; Value = "KickToken" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick" }
    @{ Cell = "D42"; Value = "'0.003815" }
    @{ Cell = "E42"; Value = "'-33.03%" }
    @{ Cell = "E43"; Value = "'-5.09%" }
    @{ Cell = "D44"; Value = "'0.009697" }
    @{ Cell = "E44"; Value = "'-2.51%" }
    @{ Cell = "D45"; Value = "'0.00005234" }
    @{ Cell = "E45"; Value = "'-4.38%" }
    @{ Cell = "E46"; Value = "'0.10%" }
    @{ Cell = "D47"; Value = "'0.09302" }
    @{ Cell = "E47"; Value = "'31.13%" }
    @{ Cell = "E48"; Value = "'-12.76%" }
    @{ Cell = "D49"; Value = "'0.00002100" }
    @{ Cell = "E49"; Value = "'0.11%" }
    @{ Cell = "D50"; Value = "'0.0002000" }
    @{ Cell = "E50"; Value = "'0.10%" }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

